{"js": "const body = context.document.body;\nconst lastParagraph = body.paragraphs.getLast();\nconst range = lastParagraph.getRange(\"End\");\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n<w:p/>\n<w:p/>\n<w:p/>\n<w:p>\n  <w:r>\n    <w:t xml:space=\"preserve\">Ho aggiornato tutti i vari commenti...togliendo quelli apposto, tenendo quelli da fare ancora e aggiungendo errori o \"problemi\" del codice ... Sono tutti in fondo alla Plancia e </w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r>\n    <w:t>Main</w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n</w:p>\n<w:p/>\n<w:p>\n  <w:r>\n    <w:t>Poi ho dato una sistemata alle frasi che appaiono in console per rendere pi\u00f9 semplice e chiaro come mettere in input (cercando cos\u00ec di evitare anche alcuni problemi con immissioni strane), si potrebbero aggiungere dei cicli di controllo al massimo</w:t>\n  </w:r>\n</w:p>\n<w:p/>\n<w:p>\n  <w:r>\n    <w:t>Ho aggiunto dei cicli di controllo quando c'\u00e8 l'immissione dei nomi per evitare duplicati e nomi vuoti</w:t>\n  </w:r>\n</w:p>\n<w:p/>\n<w:p>\n  <w:r>\n    <w:t xml:space=\"preserve\">E nella classe </w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r>\n    <w:t>cartaObiettivopersonale</w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r>\n    <w:t xml:space=\"preserve\"> ti ho creato le matrici di tutti gli obiettivi personali, cercando di velocizzare il lavoro a te</w:t>\n  </w:r>\n</w:p>\n<w:p/>\n<w:p>\n  <w:r>\n    <w:t xml:space=\"preserve\">Una cosa importante </w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r>\n    <w:t>\u00e9</w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r>\n    <w:t xml:space=\"preserve\"> che forse c'\u00e8 un problema con le iniziali di colori ( due colori iniziano con la lettera B e credo che quando hai creato le tessere per la plancia non ti sei accorto di questa cosa)... A me </w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r>\n    <w:t>\u00e9</w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r>\n    <w:t xml:space=\"preserve\"> venuto in mente di fare diventare la B di bianco una W (come in inglese) e tenere la B di blu, e cos\u00ec ho fatto per gli obiettivi personali</w:t>\n  </w:r>\n</w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nrange.insertOoxml(ooxml, \"After\");\nawait context.sync();\n", "ps1": "# Append the \"riassunto delle ultime modifiche\" paragraphs (and 3 leading\n# blank paragraphs) right after the last existing paragraph in the body,\n# just before the sectPr. Each entry below is the inner OOXML (runs /\n# proofErr marks) of one new paragraph; an empty string produces a bare\n# <w:p/>.\n$d = $word.ActiveDocument\n$ns = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'\n\n$paraInnerXmls = @(\n    '',\n    '',\n    '',\n    '<w:r><w:t xml:space=\"preserve\">Ho aggiornato tutti i vari commenti...togliendo quelli apposto, tenendo quelli da fare ancora e aggiungendo errori o \"problemi\" del codice ... Sono tutti in fondo alla Plancia e </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Main</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>',\n    '',\n    '<w:r><w:t>Poi ho dato una sistemata alle frasi che appaiono in console per rendere pi\u00f9 semplice e chiaro come mettere in input (cercando cos\u00ec di evitare anche alcuni problemi con immissioni strane), si potrebbero aggiungere dei cicli di controllo al massimo</w:t></w:r>',\n    '',\n    '<w:r><w:t>Ho aggiunto dei cicli di controllo quando c''\u00e8 l''immissione dei nomi per evitare duplicati e nomi vuoti</w:t></w:r>',\n    '',\n    '<w:r><w:t xml:space=\"preserve\">E nella classe </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>cartaObiettivopersonale</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> ti ho creato le matrici di tutti gli obiettivi personali, cercando di velocizzare il lavoro a te</w:t></w:r>',\n    '',\n    '<w:r><w:t xml:space=\"preserve\">Una cosa importante </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>\u00e9</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> che forse c''\u00e8 un problema con le iniziali di colori ( due colori iniziano con la lettera B e credo che quando hai creato le tessere per la plancia non ti sei accorto di questa cosa)... A me </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>\u00e9</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> venuto in mente di fare diventare la B di bianco una W (come in inglese) e tenere la B di blu, e cos\u00ec ho fatto per gli obiettivi personali</w:t></w:r>'\n)\n\nforeach ($inner in $paraInnerXmls) {\n    $lastPara = $d.Paragraphs.Last\n    $lastPara.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Last\n    $xml = \"<w:p xmlns:w='$ns'>$inner</w:p>\"\n    $newPara.Range.InsertXML($xml)\n}\n\nWrite-Output \"done; paragraphs=$($d.Paragraphs.Count)\""}
